$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "2008-01"
$ws.Range("A3").Value = "2008-02"
$ws.Range("A4").Value = "2008-03"
$ws.Range("A5").Value = "2008-04"
$ws.Range("A6").Value = "2008-05"
$ws.Range("A7").Value = "2008-06"
$ws.Range("A8").Value = "2008-07"
$ws.Range("A9").Value = "2008-08"
$ws.Range("A10").Value = "2008-09"
$ws.Range("A14").Value = "2009-01"
$ws.Range("A15").Value = "2009-02"
$ws.Range("A16").Value = "2009-03"
$ws.Range("A17").Value = "2009-04"
$ws.Range("A18").Value = "2009-05"
$ws.Range("A19").Value = "2009-06"
$ws.Range("A20").Value = "2009-07"
$ws.Range("A21").Value = "2009-08"
$ws.Range("A22").Value = "2009-09"
$ws.Range("A26").Value = "2010-01"
$ws.Range("A27").Value = "2010-02"
$ws.Range("A28").Value = "2010-03"
$ws.Range("A29").Value = "2010-04"
$ws.Range("A30").Value = "2010-05"
$ws.Range("A31").Value = "2010-06"
$ws.Range("A32").Value = "2010-07"
$ws.Range("A33").Value = "2010-08"
$ws.Range("A34").Value = "2010-09"
$ws.Range("A38").Value = "2011-01"
$ws.Range("A39").Value = "2011-02"
$ws.Range("A40").Value = "2011-03"
$ws.Range("A41").Value = "2011-04"
$ws.Range("A42").Value = "2011-05"
$ws.Range("A43").Value = "2011-06"
$ws.Range("A44").Value = "2011-07"
$ws.Range("A45").Value = "2011-08"
$ws.Range("A46").Value = "2011-09"
$ws.Range("A50").Value = "2012-01"
$ws.Range("A51").Value = "2012-02"
$ws.Range("A52").Value = "2012-03"
$ws.Range("A53").Value = "2012-04"
$ws.Range("A54").Value = "2012-05"
$ws.Range("A55").Value = "2012-06"
$ws.Range("A56").Value = "2012-07"
$ws.Range("A57").Value = "2012-08"
$ws.Range("A58").Value = "2012-09"
$ws.Range("A62").Value = "2013-01"
$ws.Range("A63").Value = "2013-02"
$ws.Range("A64").Value = "2013-03"
$ws.Range("A65").Value = "2013-04"
$ws.Range("A66").Value = "2013-05"
$ws.Range("A67").Value = "2013-06"
$ws.Range("A68").Value = "2013-07"
$ws.Range("A69").Value = "2013-08"
$ws.Range("A70").Value = "2013-09"
$ws.Range("A74").Value = "2014-01"
$ws.Range("A75").Value = "2014-02"
$ws.Range("A76").Value = "2014-03"
$ws.Range("A77").Value = "2014-04"
$ws.Range("A78").Value = "2014-05"
$ws.Range("A79").Value = "2014-06"
$ws.Range("A80").Value = "2014-07"
$ws.Range("A81").Value = "2014-08"
$ws.Range("A82").Value = "2014-09"
$ws.Range("A86").Value = "2015-01"
$ws.Range("A87").Value = "2015-02"
$ws.Range("A88").Value = "2015-03"
$ws.Range("A89").Value = "2015-04"
$ws.Range("A90").Value = "2015-05"
$ws.Range("A91").Value = "2015-06"
$ws.Range("A92").Value = "2015-07"
$ws.Range("A93").Value = "2015-08"
$ws.Range("A94").Value = "2015-09"
$ws.Range("A98").Value = "2016-01"
$ws.Range("A99").Value = "2016-02"
$ws.Range("A100").Value = "2016-03"
$ws.Range("A101").Value = "2016-04"
$ws.Range("A102").Value = "2016-05"
$ws.Range("A103").Value = "2016-06"
$ws.Range("A104").Value = "2016-07"
$ws.Range("A105").Value = "2016-08"
$ws.Range("A106").Value = "2016-09"
$ws.Range("A110").Value = "2017-01"
$ws.Range("A111").Value = "2017-02"
$ws.Range("A112").Value = "2017-03"
$ws.Range("A113").Value = "2017-04"
$ws.Range("A114").Value = "2017-05"
$ws.Range("A115").Value = "2017-06"
$ws.Range("A116").Value = "2017-07"
$ws.Range("A117").Value = "2017-08"
$ws.Range("A118").Value = "2017-09"
$ws.Range("A122").Value = "2018-01"
$ws.Range("A123").Value = "2018-02"
$ws.Range("A124").Value = "2018-03"
$ws.Range("A125").Value = "2018-04"
$ws.Range("A126").Value = "2018-05"
$ws.Range("A127").Value = "2018-06"
